$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create row 10 by copying the formatting/structure of row 9 (A9:T9 -> A10:T10),
# so the date cell (D10) keeps the correct date number format and no stray
# full-row cells are introduced.
$ws.Range("A9:T9").Copy($ws.Range("A10:T10"))

# Final data values for rows 2-10 (columns D, M, N, O, P, Q, S change per the diff;
# other columns are unchanged but re-asserted here for safety/clarity).

# Row 2
$ws.Range("D2").Value = 44216
$ws.Range("M2").Value = 55
$ws.Range("N2").Value = 11000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 11545
$ws.Range("Q2").Value = "$/caja 14 kilos empedrada"
$ws.Range("S2").Value = 825

# Row 3
$ws.Range("D3").Value = 44253
$ws.Range("M3").Value = 90
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 13000
$ws.Range("P3").Value = 12667
$ws.Range("Q3").Value = "$/caja 14 kilos empedrada"
$ws.Range("S3").Value = 905

# Row 4
$ws.Range("D4").Value = 44181
$ws.Range("M4").Value = 65
$ws.Range("N4").Value = 9000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 9462
$ws.Range("Q4").Value = "$/caja 14 kilos empedrada"
$ws.Range("S4").Value = 676

# Row 5
$ws.Range("D5").Value = 44172
$ws.Range("M5").Value = 90
$ws.Range("N5").Value = 8500
$ws.Range("O5").Value = 9000
$ws.Range("P5").Value = 8806
$ws.Range("Q5").Value = "$/caja 14 kilos empedrada"
$ws.Range("S5").Value = 629

# Row 6
$ws.Range("D6").Value = 44210
$ws.Range("M6").Value = 70
$ws.Range("N6").Value = 10000
$ws.Range("O6").Value = 11000
$ws.Range("P6").Value = 10357
$ws.Range("Q6").Value = "$/caja 14 kilos empedrada"
$ws.Range("S6").Value = 740

# Row 7
$ws.Range("D7").Value = 44232
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = 11000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 11583
$ws.Range("Q7").Value = "$/caja 14 kilos empedrada"
$ws.Range("S7").Value = 827

# Row 8
$ws.Range("D8").Value = 44229
$ws.Range("M8").Value = 55
$ws.Range("N8").Value = 11000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 11364
$ws.Range("Q8").Value = "$/caja 14 kilos empedrada"
$ws.Range("S8").Value = 812

# Row 9
$ws.Range("D9").Value = 45138
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 14000
$ws.Range("P9").Value = 14000
$ws.Range("Q9").Value = "$/caja 14 kilos granel"
$ws.Range("S9").Value = 1000

# Row 10 (new row)
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").Value = 45140
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100102
$ws.Range("H10").Value = "Cítricos"
$ws.Range("I10").Value = 100102006
$ws.Range("J10").Value = "Pomelo"
$ws.Range("K10").Value = "Start Ruby"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 30
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 15000
$ws.Range("Q10").Value = "$/caja 14 kilos granel"
$ws.Range("R10").Value = "Región de O'Higgins"
$ws.Range("S10").Value = 1071
$ws.Range("T10").Value = 14
